$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update grade values (column C) for several subjects; column D (Jegy*kr)
# recalculates automatically via its shared formula B*C.
$ws.Range("C4").Value = 3
$ws.Range("C7").Value = 5
$ws.Range("C9").Value = 3
$ws.Range("C12").Value = 2

# Highlight the updated cells (C7:C8) with a green fill.
$ws.Range("C7:C8").Interior.Color = 0x50B000

# Move the active selection to C14.
$ws.Range("C14").Select()
